$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: G19 gets a value of 5 (style s="2" stays unchanged) ---
$ws.Range("G19").Value = 5

# --- Row 31: G31 gets a value of 5 (style s="2" stays unchanged) ---
$ws.Range("G31").Value = 5

# --- Row 32 ---
# F32 changes fill/format from the green "answer" style (s="5") to the
# plain style (s="2") already used by G32/H32/I32. Copy that format from
# a neighbouring cell that already carries style s="2" (G32) so the
# shared style index is reused instead of minting a new one.
$ws.Range("G32").Copy() | Out-Null
$ws.Range("F32").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("F32").Value = 5

# G32 gets a value of 5 (style s="2" stays unchanged)
$ws.Range("G32").Value = 5

# --- Selection: move the active cell of the bottom-right frozen pane to G32 ---
$ws.Range("G32").Select() | Out-Null
